$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1516.12
$ws.Range("I33").Value = 988.8333
$ws.Range("K33").Value = 988.8333
$ws.Range("M33").Value = -759.8333

$ws.Range("H98").Value = 18726828
$ws.Range("I98").Value = 6897871.5
$ws.Range("J98").Value = 56842356
$ws.Range("K98").Value = 6897871.5
$ws.Range("L98").Value = 56842356
$ws.Range("M98").Value = -6896373.5
$ws.Range("N98").Value = -56845352

$ws.Range("H107").Value = 8000
$ws.Range("I107").Value = 8000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 8000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -6080
$ws.Range("N107").ClearContents()

$ws.Range("H116").Value = 5101007.5
$ws.Range("I116").Value = 11113877
$ws.Range("J116").Value = 4199077
$ws.Range("K116").Value = 11113877
$ws.Range("L116").Value = 4199077
$ws.Range("M116").Value = -11110435
$ws.Range("N116").Value = -4205961

$ws.Range("H122").Value = 18726828
$ws.Range("I122").Value = 6897871.5
$ws.Range("J122").Value = 56842356
$ws.Range("K122").Value = 20693614.5
$ws.Range("L122").Value = 170527068
$ws.Range("M122").Value = -20691164.5
$ws.Range("N122").Value = -170531968

$ws.Range("H128").Value = 38828.668
$ws.Range("J128").Value = 38828.668
$ws.Range("L128").Value = 38828.668
$ws.Range("N128").Value = -48788.668

$ws.Range("H132").Value = 2527259
$ws.Range("I132").Value = 1860.1842
$ws.Range("J132").Value = 18521452
$ws.Range("K132").Value = 5580.5526
$ws.Range("L132").Value = 55564356
$ws.Range("M132").Value = -3050.5526
$ws.Range("N132").Value = -55569416

$ws.Range("H137").Value = 48155850
$ws.Range("I137").Value = 125000760
$ws.Range("J137").Value = 22540880
$ws.Range("K137").Value = 375002280
$ws.Range("L137").Value = 67622640
$ws.Range("M137").Value = -374999730
$ws.Range("N137").Value = -67627740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8778050
$ws.Range("I32").Value = 6260.041
$ws.Range("J32").Value = 62505264
$ws.Range("K32").Value = 6260.041
$ws.Range("L32").Value = 62505264
$ws.Range("M32").Value = -5973.041
$ws.Range("N32").Value = -62505838

$ws.Range("H61").Value = 4293915
$ws.Range("I61").Value = 3631256.5
$ws.Range("J61").Value = 5884295
$ws.Range("K61").Value = 3631256.5
$ws.Range("L61").Value = 5884295
$ws.Range("M61").Value = -3631044.5
$ws.Range("N61").Value = -5884719

$ws.Range("H74").Value = 68481260
$ws.Range("I74").Value = 115742120
$ws.Range("J74").Value = 38099292
$ws.Range("K74").Value = 115742120
$ws.Range("L74").Value = 38099292
$ws.Range("M74").Value = -115741246
$ws.Range("N74").Value = -38101040

$ws.Range("H77").Value = 68481260
$ws.Range("I77").Value = 115742120
$ws.Range("J77").Value = 38099292
$ws.Range("K77").Value = 578710600
$ws.Range("L77").Value = 190496460
$ws.Range("M77").Value = -578706232
$ws.Range("N77").Value = -190505196

$ws.Range("H88").Value = 6944.4443
$ws.Range("I88").Value = 2166.6667
$ws.Range("K88").Value = 2166.6667
$ws.Range("M88").Value = -1760.6667

$ws.Range("H91").Value = 6944.4443
$ws.Range("I91").Value = 2166.6667
$ws.Range("K91").Value = 2166.6667
$ws.Range("M91").Value = -762.6667000000002

$ws.Range("H132").Value = 23723230
$ws.Range("I132").Value = 24081548
$ws.Range("J132").Value = 23234614
$ws.Range("K132").Value = 72244644
$ws.Range("L132").Value = 69703842
$ws.Range("M132").Value = -72242114
$ws.Range("N132").Value = -69708902

$ws.Range("H136").Value = 4293915
$ws.Range("I136").Value = 3631256.5
$ws.Range("J136").Value = 5884295
$ws.Range("K136").Value = 10893769.5
$ws.Range("L136").Value = 17652885
$ws.Range("M136").Value = -10891219.5
$ws.Range("N136").Value = -17657985

$ws.Range("H138").Value = 59644.25
$ws.Range("J138").Value = 59644.25
$ws.Range("L138").Value = 59644.25
$ws.Range("N138").Value = -69924.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22059780
$ws.Range("I134").Value = 31250844
$ws.Range("J134").Value = 3677649
$ws.Range("K134").Value = 93752532
$ws.Range("L134").Value = 11032947
$ws.Range("M134").Value = -93749997
$ws.Range("N134").Value = -11038017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12347124
$ws.Range("I31").Value = 22728302
$ws.Range("J31").Value = 1939
$ws.Range("K31").Value = 22728302
$ws.Range("L31").Value = 1939
$ws.Range("M31").Value = -22728007
$ws.Range("N31").Value = -2529

$ws.Range("H34").Value = 12347124
$ws.Range("I34").Value = 22728302
$ws.Range("J34").Value = 1939
$ws.Range("K34").Value = 22728302
$ws.Range("L34").Value = 1939
$ws.Range("M34").Value = -22728100
$ws.Range("N34").Value = -2343

$ws.Range("H105").Value = 20833.334
$ws.Range("I105").Value = 10000
$ws.Range("J105").Value = 26250
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 26250
$ws.Range("M105").Value = -8253
$ws.Range("N105").Value = -29744

$ws.Range("H107").Value = 486.27585
$ws.Range("I107").Value = 293.83334
$ws.Range("J107").Value = 801.1818
$ws.Range("K107").Value = 293.83334
$ws.Range("L107").Value = 801.1818
$ws.Range("M107").Value = 1626.16666
$ws.Range("N107").Value = -4641.1818

$ws.Range("H132").Value = 1726046
$ws.Range("I132").Value = 2778750.5
$ws.Range("J132").Value = 3438.7273
$ws.Range("K132").Value = 8336251.5
$ws.Range("L132").Value = 10316.1819
$ws.Range("M132").Value = -8333721.5
$ws.Range("N132").Value = -15376.1819

$ws.Range("H134").Value = 1144875.2
$ws.Range("I134").Value = 1524.7587
$ws.Range("J134").Value = 6671069
$ws.Range("K134").Value = 4574.2761
$ws.Range("L134").Value = 20013207
$ws.Range("M134").Value = -2039.2761
$ws.Range("N134").Value = -20018277

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 229.6
$ws.Range("J52").Value = 229.6
$ws.Range("L52").Value = 688.8
$ws.Range("N52").Value = -1220.8

$ws.Range("H131").Value = 7368545.5
$ws.Range("I131").Value = 83333544
$ws.Range("J131").Value = 17093.791
$ws.Range("K131").Value = 250000632
$ws.Range("L131").Value = 51281.37300000001
$ws.Range("M131").Value = -249995592
$ws.Range("N131").Value = -61361.37300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 21841.084
$ws.Range("I113").Value = 998.75
$ws.Range("J113").Value = 63525.75
$ws.Range("K113").Value = 998.75
$ws.Range("L113").Value = 63525.75
$ws.Range("M113").Value = 1171.25
$ws.Range("N113").Value = -67865.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2454.3333
$ws.Range("I7").Value = 2412.7144
$ws.Range("K7").Value = 2412.7144
$ws.Range("M7").Value = -2300.7144

$ws.Range("H40").Value = 3425.4546
$ws.Range("I40").Value = 3355.5557
$ws.Range("J40").Value = 3740
$ws.Range("K40").Value = 3355.5557
$ws.Range("L40").Value = 3740
$ws.Range("M40").Value = -3219.5557
$ws.Range("N40").Value = -4012

$ws.Range("H61").Value = 2750.818
$ws.Range("I61").Value = 2181.75
$ws.Range("J61").Value = 4268.3335
$ws.Range("K61").Value = 2181.75
$ws.Range("L61").Value = 4268.3335
$ws.Range("M61").Value = -1979.75
$ws.Range("N61").Value = -4672.3335

$ws.Range("H113").Value = 2750.818
$ws.Range("I113").Value = 2181.75
$ws.Range("J113").Value = 4268.3335
$ws.Range("K113").Value = 2181.75
$ws.Range("L113").Value = 4268.3335
$ws.Range("M113").Value = -11.75
$ws.Range("N113").Value = -8608.333500000001

$ws.Range("H126").Value = 2454.3333
$ws.Range("I126").Value = 2412.7144
$ws.Range("K126").Value = 7238.1432
$ws.Range("M126").Value = -4768.1432

$ws.Range("H132").Value = 4532710.5
$ws.Range("I132").Value = 6419881.5
$ws.Range("K132").Value = 19259644.5
$ws.Range("M132").Value = -19257114.5

$ws.Range("H136").Value = 10620904
$ws.Range("I136").Value = 3274446.5
$ws.Range("K136").Value = 9823339.5
$ws.Range("M136").Value = -9820789.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 416
$ws.Range("I113").Value = 430.25
$ws.Range("J113").Value = 302
$ws.Range("K113").Value = 1290.75
$ws.Range("L113").Value = 906
$ws.Range("M113").Value = 879.25
$ws.Range("N113").Value = -5246

$ws.Range("H122").Value = 1576.8572
$ws.Range("I122").Value = 1454.4445
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 4363.333500000001
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -1913.333500000001
$ws.Range("N122").Value = -10870

$ws.Range("H132").Value = 336005.94
$ws.Range("I132").Value = 454723.16
$ws.Range("K132").Value = 1364169.48
$ws.Range("M132").Value = -1361639.48

$ws.Range("H136").Value = 5293.595
$ws.Range("I136").Value = 3962.9583
$ws.Range("K136").Value = 11888.8749
$ws.Range("M136").Value = -9338.874899999999
